# This workbook ("AHB-Diff") lists, side by side, the segments/fields of an
# "old" and a "new" EDIFACT message format definition. This edit renames the
# generic "_old" / "_new" header suffixes to the concrete format versions
# being compared (FV2210 and FV2304), turns the header row + data range into
# a proper Excel Table (so the headers carry AutoFilter drop-downs), and
# freezes the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row -------------------------------------------------
# Columns A-J describe the "old" format -> suffix becomes "_FV2210"
foreach ($col in 1..10) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2.Replace("_old", "_FV2210")
}

# Column K is just "diff" and is left untouched.

# Columns L-U describe the "new" format -> suffix becomes "_FV2304"
foreach ($col in 12..21) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2.Replace("_new", "_FV2304")
}

# --- 2) Turn the data range into an Excel Table ---------------------------
$dataRange = $ws.Range("A1:U70")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- 3) Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
